$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") rows 2-37, replacing old Strike# values
$values = @{
    2  = 3
    3  = 0
    4  = 7
    5  = 6
    6  = 6
    7  = 4
    8  = 5
    9  = 1
    10 = 3
    11 = 0
    12 = 3
    13 = 6
    14 = 7
    15 = 4
    16 = 4
    17 = 3
    18 = 4
    19 = 7
    20 = 6
    21 = 6
    22 = 5
    23 = 4
    24 = 5
    25 = 3
    26 = 4
    27 = 4
    28 = 5
    29 = 5
    30 = 3
    31 = 9
    32 = 3
    33 = 2
    34 = 3
    35 = 1
    36 = 2
    37 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
